$d = $word.ActiveDocument

function New-RunsPackageXml([string[]]$texts, [string[]]$langs) {
    $w_ns = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
    $runsXml = ""
    for ($i = 0; $i -lt $texts.Length; $i++) {
        $t = $texts[$i]
        $lang = $langs[$i]
        if ($t -ne $t.Trim()) {
            $preserve = " xml:space=`"preserve`""
        } else {
            $preserve = ""
        }
        $runsXml += "<w:r><w:rPr><w:rFonts w:ascii=`"Proxima Nova`" w:eastAsia=`"Proxima Nova`" w:hAnsi=`"Proxima Nova`" w:cs=`"Proxima Nova`"/><w:lang w:val=`"$lang`"/></w:rPr><w:t$preserve>$t</w:t></w:r>"
    }
    $pkg = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="' + $w_ns + '"><w:body><w:p>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    return $pkg
}

# --- Edit 1: split the opening sentence into three runs ---
# The Find target is extended to also swallow the immediately-following
# "Нетологии" run (wrapped in <w:proofErr> spell-check markers). Ending
# the replaced Range exactly on the zero-width <w:proofErr> boundary is
# ambiguous, so instead we include that whole run in the replacement and
# re-emit it byte-for-byte unchanged - this keeps the proofErr markers
# anchored correctly relative to it.
$f1 = $d.Content
$f1.Find.Execute(
    "Я менеджер отдела снабжения на заводе по производству тормозного оборудования для грузовых вагонов. С октября 2021 года прохожу обучение на портале Нетологии",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $f1.Find.Found) {
    throw "Edit 1 target text not found"
}

# Re-derive a plain Range so InsertXML replaces the span in place
# instead of appending after the Find-narrowed range.
$r1 = $d.Range($f1.Start, $f1.End)

$w_ns = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$frag1 = New-RunsPackageXml @(
    "Я ",
    "специалист отдел тестирования продукта, обеспечивающего взаимодействие и обмен данными между банками и СМЭВ",
    ". С октября 2021 года прохожу обучение на портале "
) @("ru-RU", "ru-RU", "ru-RU")

# Splice the unchanged "Нетологии" run (with its spell-check markers)
# back in right before the closing </w:p>, preserving its original rPr.
$netologiiXml = '<w:proofErr w:type="spellStart"/><w:r w:rsidRPr="002B2DB9"><w:rPr><w:rFonts w:ascii="Proxima Nova" w:eastAsia="Proxima Nova" w:hAnsi="Proxima Nova" w:cs="Proxima Nova"/><w:lang w:val="ru-RU"/></w:rPr><w:t>Нетологии</w:t></w:r><w:proofErr w:type="spellEnd"/>'
$frag1 = $frag1.Replace("</w:p></w:body>", ($netologiiXml + "</w:p></w:body>"))

$r1.InsertXML($frag1)

# --- Edit 2: split the closing sentence into seven runs ---
$f2 = $d.Content
$f2.Find.Execute(
    " для тестировщиков, автоматизированное тестирование, а также сдал курсовой проект по тестированию веб-сервиса.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $f2.Find.Found) {
    throw "Edit 2 target text not found"
}

$r2 = $d.Range($f2.Start, $f2.End)

$texts2 = @(
    " для тестировщиков, ",
    "JavaScript",
    ",",
    " ",
    "автоматизированное тестирование,",
    " тестирование производительности,",
    " а также сдал курсовой проект по тестированию веб-сервиса."
)
$langs2 = @("ru-RU", "en-US", "ru-RU", "ru-RU", "ru-RU", "ru-RU", "ru-RU")

$r2.InsertXML((New-RunsPackageXml $texts2 $langs2))
